$wb = $excel.ActiveWorkbook

# --- Sheet "Productos": fix stock values ---
$wsProductos = $wb.Worksheets.Item("Productos")
$wsProductos.Cells.Item(3, 3).Value = 44
$wsProductos.Cells.Item(4, 3).Value = 67
$wsProductos.Cells.Item(5, 3).Value = 52

# --- Sheet "Ventas": add new sales rows ---
$wsVentas = $wb.Worksheets.Item("Ventas")

$wsVentas.Cells.Item(6, 1).Value = "V-376EB348"
$wsVentas.Cells.Item(6, 2).Value = "2025-08-01 20:02:27"
$wsVentas.Cells.Item(6, 3).Value = 33000

$wsVentas.Cells.Item(7, 1).Value = "V-A6E52CA1"
$wsVentas.Cells.Item(7, 2).Value = "2025-08-02 00:58:53"
$wsVentas.Cells.Item(7, 3).Value = 25000

$wsVentas.Cells.Item(8, 1).Value = "V-20BFBEEF"
$wsVentas.Cells.Item(8, 2).Value = "2025-08-02 00:59:00"
$wsVentas.Cells.Item(8, 3).Value = 10000

$wsVentas.Cells.Item(9, 1).Value = "V-686EFA3C"
$wsVentas.Cells.Item(9, 2).Value = "2025-08-02 00:59:05"
$wsVentas.Cells.Item(9, 3).Value = 60000

# --- Sheet "RegistroCaja": add new cash register entries ---
$wsCaja = $wb.Worksheets.Item("RegistroCaja")

$wsCaja.Cells.Item(5, 1).Value = "2025-08-01 20:02:27"
$wsCaja.Cells.Item(5, 2).Value = "VENTA - ID: V-376EB348 | Monto: `$33000,00"

$wsCaja.Cells.Item(6, 1).Value = "2025-08-02 00:58:53"
$wsCaja.Cells.Item(6, 2).Value = "VENTA - ID: V-A6E52CA1 | Monto: `$25000,00"

$wsCaja.Cells.Item(7, 1).Value = "2025-08-02 00:59:00"
$wsCaja.Cells.Item(7, 2).Value = "VENTA - ID: V-20BFBEEF | Monto: `$10000,00"

$wsCaja.Cells.Item(8, 1).Value = "2025-08-02 00:59:05"
$wsCaja.Cells.Item(8, 2).Value = "VENTA - ID: V-686EFA3C | Monto: `$60000,00"
